# kop en voettekst toegevoegd
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch to Page Layout view, as Excel does automatically when you add a
# header/footer via the UI (reflected by the sheetView's view="pageLayout"
# in the saved file).
$ws.Activate()
$excel.ActiveWindow.View = [Microsoft.Office.Interop.Excel.XlWindowView]::xlPageLayoutView

# Move the active selection to G3, as recorded in the saved view state.
$ws.Range("G3").Select()

# Add the centered header and the two-line centered footer.
$ps = $ws.PageSetup
$ps.CenterHeader = "Datamodel Database"
$ps.CenterFooter = "Datum: 17-4-14" + [char]10 + "Projectleden: Menno, Rick, Regilio, Sharif en Danny"
